$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E5").Value = "❌ EXPIRED 3758 days ago"
$ws.Range("E6").Value = "⚠️ Expires in 21 days"
$ws.Range("E7").Value = "⚠️ Expires in 21 days"
$ws.Range("E8").Value = "⚠️ Expires in 21 days"
